# Apply the cryptos-list refresh described in the commit:
#   "Updated cryptos list on Sun Mar  5 19:50:46 UTC 2023 with GitHub Actions"
#
# Every target cell on Sheet1 is plain text (t="inlineStr" in the original
# OOXML), including price cells that look numeric (e.g. "1.002", "0.3689").
# Excel's COM Range.Value setter auto-coerces a clean numeric-looking string
# into a real number, which would both change the stored type and silently
# drop meaningful trailing zeros (e.g. "149.00" -> 149, "0.07530" -> 0.0753).
# To keep those cells as text we prefix an apostrophe, exactly like a user
# typing '149.00 into the cell would - Excel stores the text verbatim (without
# the apostrophe) and only flags the cell with its "stored as text" marker.
# Values that already cannot parse as a number (e.g. "22.457.72", which has
# two dots) are assigned as-is since Excel leaves them as text naturally.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '22.457.72'
$ws.Range("E2").Value = '  +0.44%  '
$ws.Range("D3").Value = '1.569.59'
$ws.Range("E3").Value = '  +0.19%  '
$ws.Range("E4").Value = '  -0.15%  '
$ws.Range("D5").Value = '''1.002'
$ws.Range("D6").Value = '''290.16'
$ws.Range("E6").Value = '  -0.13%  '
$ws.Range("D7").Value = '''0.3689'
$ws.Range("E7").Value = '  -1.34%  '
$ws.Range("D8").Value = '''49.89'
$ws.Range("E8").Value = '  +1.89%  '
$ws.Range("D9").Value = '''0.3370'
$ws.Range("E9").Value = '  -0.32%  '
$ws.Range("D10").Value = '''1.144'
$ws.Range("E10").Value = '  +1.45%  '
$ws.Range("D11").Value = '''0.07530'
$ws.Range("E11").Value = '  +0.11%  '
$ws.Range("E12").Value = '  -0.16%  '
$ws.Range("D13").Value = '''21.13'
$ws.Range("E13").Value = '  +1.59%  '
$ws.Range("D14").Value = '''6.021'
$ws.Range("E14").Value = '  +1.83%  '
$ws.Range("D15").Value = '''6.970'
$ws.Range("E15").Value = '  +1.34%  '
$ws.Range("D16").Value = '1.570.63'
$ws.Range("E16").Value = '  +0.32%  '
$ws.Range("D17").Value = '''0.00001120'
$ws.Range("E17").Value = '  +0.20%  '
$ws.Range("D18").Value = '''90.39'
$ws.Range("E18").Value = '  +0.92%  '
$ws.Range("D19").Value = '''0.06771'
$ws.Range("E19").Value = '  +0.71%  '
$ws.Range("D20").Value = '''1.002'
$ws.Range("E20").Value = '  -0.14%  '
$ws.Range("D21").Value = '''6.363'
$ws.Range("E21").Value = '  +3.23%  '
$ws.Range("D22").Value = '''16.38'
$ws.Range("E22").Value = '  -0.16%  '
$ws.Range("E23").Value = '  +2.89%  '
$ws.Range("D24").Value = '22.450.45'
$ws.Range("E24").Value = '  +0.39%  '
$ws.Range("E25").Value = '  +0.13%  '
$ws.Range("D26").Value = '''2.645'
$ws.Range("E26").Value = '  -2.17%  '
$ws.Range("E27").Value = '  +0.05%  '
$ws.Range("D28").Value = '''149.00'
$ws.Range("E28").Value = '  +0.98%  '
$ws.Range("D29").Value = '''5.051'
$ws.Range("E29").Value = '  +0.25%  '
$ws.Range("D30").Value = '''124.95'
$ws.Range("E30").Value = '  +0.11%  '
$ws.Range("D31").Value = '1.750.05'
$ws.Range("E31").Value = '  +0.58%  '
$ws.Range("D32").Value = '''1.063'
$ws.Range("E32").Value = '  +8.14%  '
$ws.Range("D33").Value = '''6.181'
$ws.Range("E33").Value = '  +3.41%  '
$ws.Range("D34").Value = '''2.012'
$ws.Range("E34").Value = '  -0.30%  '
$ws.Range("D35").Value = '''9.780'
$ws.Range("E35").Value = '  -1.91%  '
$ws.Range("E36").Value = '  -1.44%  '
$ws.Range("B37").Value = 'VeChain'
$ws.Range("C37").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D37").Value = '''0.02467'
$ws.Range("E37").Value = '  -0.19%  '
$ws.Range("B38").Value = 'TrustWalletToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D38").Value = '''1.356'
$ws.Range("E38").Value = '  -3.44%  '
$ws.Range("D39").Value = '''0.2299'
$ws.Range("E39").Value = '  +1.11%  '
$ws.Range("E40").Value = '  +2.05%  '
$ws.Range("D41").Value = '''5.410'
$ws.Range("E41").Value = '  +0.83%  '
$ws.Range("D42").Value = '''11.20'
$ws.Range("E42").Value = '  +2.05%  '
$ws.Range("D43").Value = '''0.6220'
$ws.Range("E43").Value = '  -0.28%  '
$ws.Range("D44").Value = '''14.14'
$ws.Range("E44").Value = '  +1.27%  '
$ws.Range("E45").Value = '  -0.09%  '
$ws.Range("D46").Value = '''3.804'
$ws.Range("E46").Value = '  +0.20%  '
$ws.Range("D47").Value = '''0.5848'
$ws.Range("E47").Value = '  -0.30%  '
$ws.Range("D48").Value = '''2.068'
$ws.Range("E48").Value = '  +0.90%  '
$ws.Range("D49").Value = '''127.69'
$ws.Range("E49").Value = '  +2.81%  '
$ws.Range("D50").Value = '''1.239'
$ws.Range("E50").Value = '  -0.93%  '
$ws.Range("D51").Value = '''0.07304'
$ws.Range("E51").Value = '  -0.14%  '
